$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# Row 33
$ws.Range("H33").Value = 1306.7059
$ws.Range("I33").Value = 917.53845
$ws.Range("J33").Value = 2571.5
$ws.Range("K33").Value = 917.53845
$ws.Range("L33").Value = 2571.5
$ws.Range("M33").Value = -688.53845
$ws.Range("N33").Value = -3029.5

# Row 98
$ws.Range("H98").Value = 155.5
$ws.Range("I98").Value = 155.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 155.5
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 1342.5

# Row 122
$ws.Range("H122").Value = 155.5
$ws.Range("I122").Value = 155.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 466.5
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 1983.5

# Row 132
$ws.Range("H132").Value = 5258.381
$ws.Range("I132").Value = 5201.647
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 15604.941
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -13074.941
$ws.Range("N132").Value = -21558.5

# Row 138
$ws.Range("H138").Value = 1862.091
$ws.Range("I138").Value = 1447.258
$ws.Range("J138").Value = 2141.652
$ws.Range("K138").Value = 4341.774
$ws.Range("L138").Value = 6424.956
$ws.Range("M138").Value = 798.2259999999997
$ws.Range("N138").Value = -16704.956

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 411263.16
$ws.Range("I32").Value = 455619.72
$ws.Range("J32").Value = 21911
$ws.Range("K32").Value = 455619.72
$ws.Range("L32").Value = 21911
$ws.Range("M32").Value = -455332.72
$ws.Range("N32").Value = -22485

# Row 74
$ws.Range("H74").Value = 912.4583
$ws.Range("I74").Value = 820.3684
$ws.Range("J74").Value = 1262.4
$ws.Range("K74").Value = 820.3684
$ws.Range("L74").Value = 1262.4
$ws.Range("M74").Value = 53.63160000000005
$ws.Range("N74").Value = -3010.4

# Row 77
$ws.Range("H77").Value = 912.4583
$ws.Range("I77").Value = 820.3684
$ws.Range("J77").Value = 1262.4
$ws.Range("K77").Value = 4101.842
$ws.Range("L77").Value = 6312
$ws.Range("M77").Value = 266.1580000000004
$ws.Range("N77").Value = -15048

# Row 93
$ws.Range("H93").Value = 80000
$ws.Range("J93").Value = 80000
$ws.Range("L93").Value = 80000
$ws.Range("N93").Value = -84992

# Row 95
$ws.Range("H95").Value = 100000
$ws.Range("J95").Value = 100000
$ws.Range("L95").Value = 100000
$ws.Range("N95").Value = -105492

# Row 102
$ws.Range("H102").Value = 2549.3125
$ws.Range("I102").Value = 2579.2666
$ws.Range("K102").Value = 2579.2666
$ws.Range("M102").Value = -957.2665999999999

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# Row 134
$ws.Range("H134").Value = 2403.8518
$ws.Range("I134").Value = 1994.3529
$ws.Range("J134").Value = 3100
$ws.Range("K134").Value = 5983.0587
$ws.Range("L134").Value = 9300
$ws.Range("M134").Value = -3448.0587
$ws.Range("N134").Value = -14370

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# Row 132
$ws.Range("H132").Value = 4506427.5
$ws.Range("I132").Value = 1796.8948
$ws.Range("J132").Value = 9261315
$ws.Range("K132").Value = 5390.6844
$ws.Range("L132").Value = 27783945
$ws.Range("M132").Value = -2860.6844
$ws.Range("N132").Value = -27789005

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# Row 12
$ws.Range("H12").Value = 83.07692
$ws.Range("I12").Value = 4
$ws.Range("J12").Value = 97.454544
$ws.Range("K12").Value = 12
$ws.Range("L12").Value = 292.363632
$ws.Range("M12").Value = 161
$ws.Range("N12").Value = -638.3636320000001

# Row 109
$ws.Range("H109").Value = 3822.1428
$ws.Range("I109").Value = 1383.3334
$ws.Range("J109").Value = 5651.25
$ws.Range("K109").Value = 4150.0002
$ws.Range("L109").Value = 16953.75
$ws.Range("M109").Value = -3110.0002
$ws.Range("N109").Value = -19033.75

# Row 131
$ws.Range("H131").Value = 1040.6666
$ws.Range("J131").Value = 1122.6487
$ws.Range("L131").Value = 3367.9461
$ws.Range("N131").Value = -13447.9461

# Row 132
$ws.Range("H132").Value = 1718.8695
$ws.Range("I132").Value = 1257.5
$ws.Range("J132").Value = 1964.9333
$ws.Range("K132").Value = 11317.5
$ws.Range("L132").Value = 17684.3997
$ws.Range("M132").Value = -8787.5
$ws.Range("N132").Value = -22744.3997

# Row 134
$ws.Range("H134").Value = 6086.737
$ws.Range("I134").Value = 2883.0833
$ws.Range("J134").Value = 7565.346
$ws.Range("K134").Value = 8649.249899999999
$ws.Range("L134").Value = 22696.038
$ws.Range("M134").Value = -3579.249899999999
$ws.Range("N134").Value = -32836.038

# Row 139
$ws.Range("H139").Value = 3700.342
$ws.Range("I139").Value = 1077.5
$ws.Range("K139").Value = 3232.5
$ws.Range("M139").Value = 1907.5

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# Row 80
$ws.Range("H80").Value = 36360670
$ws.Range("I80").Value = 50903744
$ws.Range("J80").Value = 2997.5
$ws.Range("K80").Value = 50903744
$ws.Range("L80").Value = 2997.5
$ws.Range("M80").Value = -50902746
$ws.Range("N80").Value = -4993.5

# Row 83
$ws.Range("H83").Value = 36360670
$ws.Range("I83").Value = 50903744
$ws.Range("J83").Value = 2997.5
$ws.Range("K83").Value = 254518720
$ws.Range("L83").Value = 14987.5
$ws.Range("M83").Value = -254513728
$ws.Range("N83").Value = -24971.5

# Row 97
$ws.Range("H97").Value = 1585
$ws.Range("I97").Value = 1332.8572
$ws.Range("J97").Value = 1938
$ws.Range("K97").Value = 1332.8572
$ws.Range("L97").Value = 1938
$ws.Range("M97").Value = -836.8571999999999
$ws.Range("N97").Value = -2930

# Row 122
$ws.Range("H122").Value = 1580.4546
$ws.Range("I122").Value = 1653.7222
$ws.Range("K122").Value = 4961.1666
$ws.Range("M122").Value = -2511.1666

# Row 132
$ws.Range("H132").Value = 2527.8206
$ws.Range("I132").Value = 2187.5
$ws.Range("J132").Value = 2968.2354
$ws.Range("K132").Value = 6562.5
$ws.Range("L132").Value = 8904.706200000001
$ws.Range("M132").Value = -4032.5
$ws.Range("N132").Value = -13964.7062

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Range("H7").Value = 3247.238
$ws.Range("I7").Value = 3139.3142
$ws.Range("J7").Value = 3786.8572
$ws.Range("K7").Value = 3139.3142
$ws.Range("L7").Value = 3786.8572
$ws.Range("M7").Value = -3027.3142
$ws.Range("N7").Value = -4010.8572

# Row 40
$ws.Range("H40").Value = 118389.22
$ws.Range("I40").Value = 174500.5
$ws.Range("J40").Value = 6166.6665
$ws.Range("K40").Value = 174500.5
$ws.Range("L40").Value = 6166.6665
$ws.Range("M40").Value = -174364.5
$ws.Range("N40").Value = -6438.6665

# Row 55
$ws.Range("H55").Value = 645.53845
$ws.Range("I55").Value = 453.46667
$ws.Range("J55").Value = 907.4545000000001
$ws.Range("K55").Value = 453.46667
$ws.Range("L55").Value = 907.4545000000001
$ws.Range("M55").Value = -280.46667
$ws.Range("N55").Value = -1253.4545

# Row 93
$ws.Range("H93").Value = 10101.583
$ws.Range("I93").Value = 10855.909
$ws.Range("J93").Value = 1804
$ws.Range("K93").Value = 10855.909
$ws.Range("L93").Value = 1804
$ws.Range("M93").Value = -9607.909
$ws.Range("N93").Value = -4300

# Row 126
$ws.Range("H126").Value = 3247.238
$ws.Range("I126").Value = 3139.3142
$ws.Range("J126").Value = 3786.8572
$ws.Range("K126").Value = 9417.942599999998
$ws.Range("L126").Value = 11360.5716
$ws.Range("M126").Value = -6947.942599999998
$ws.Range("N126").Value = -16300.5716

# Row 136
$ws.Range("H136").Value = 13891085
$ws.Range("I136").Value = 1981.25
$ws.Range("J136").Value = 41669292
$ws.Range("K136").Value = 5943.75
$ws.Range("L136").Value = 125007876
$ws.Range("M136").Value = -3393.75
$ws.Range("N136").Value = -125012976
